# peakList accessors for MSPeakListsSet
# Adds "X" marks to the feature-support matrix on the "mslists" sheet
# for the newly supported peakList-related accessors.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mslists")

# Mark new cells with "X" (same shared string used elsewhere in the sheet)
$ws.Range("F7").Value = "X"
$ws.Range("G7").Value = "X"

$ws.Range("G10").Value = "X"

$ws.Range("G11").Value = "X"

$ws.Range("G12").Value = "X"

$ws.Range("F13").Value = "X"
$ws.Range("G13").Value = "X"

# Update the active selection to reflect where the user ended up (row 15)
$ws.Activate()
$ws.Range("G15").Select()
